$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.432.94'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.25%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.180.07'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.16%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.27'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.59%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.08'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.94%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.185.33'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.34%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.87%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.29%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.20%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000269'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.41%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.89'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.694.81'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.235.41'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.26'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.38%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.180.25'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.19%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '508.52'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.86%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.82'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +5.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.726'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.57%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.20'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -5.73%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.89'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.18%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.00'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.01%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.19'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.99'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.59%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.52%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.92'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.58%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.07%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.11%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.57'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.13'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.94%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0909'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '478.02'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.35%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0422'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.25%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -5.52%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.88'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.57%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.006.32'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.34%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.38%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.289'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.64%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.46'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.37%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₃0618'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +5.20%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.76'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.65%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.01%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.45%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.96%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '119.70'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.18%  '
